$wb = $excel.ActiveWorkbook

# --- Summary sheet: Total Trades 7 -> 8, Win Rate % 42.86 -> 37.5 ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 8
$summary.Range("B9").Value = 37.5

# --- Strategy Status sheet: MarketMaking Trades 7 -> 8, Win Rate % 42.86 -> 37.5 ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 8
$status.Range("G4").Value = 37.5

# --- New trade #8 row, appended to both "All Trades" and "MarketMaking" sheets ---
$newRow = @{
    A = 8
    B = "2026-02-17"
    C = "12:27:14"
    D = "MarketMaking"
    E = "DOWN"
    F = 0.88
    G = 0.88
    H = "CLOSED"
    I = 0
    J = 0
    K = 100.01
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.14
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A9").Value = $newRow.A
    # Date-like text must be forced to text format first, otherwise Excel
    # auto-converts "2026-02-17" into a date serial number.
    $ws.Range("B9").NumberFormat = "@"
    $ws.Range("B9").Value = $newRow.B
    $ws.Range("C9").Value = $newRow.C
    $ws.Range("D9").Value = $newRow.D
    $ws.Range("E9").Value = $newRow.E
    $ws.Range("F9").Value = $newRow.F
    $ws.Range("G9").Value = $newRow.G
    $ws.Range("H9").Value = $newRow.H
    $ws.Range("I9").Value = $newRow.I
    $ws.Range("J9").Value = $newRow.J
    $ws.Range("K9").Value = $newRow.K
    $ws.Range("L9").Value = $newRow.L
    $ws.Range("M9").Value = $newRow.M
    $ws.Range("N9").Value = $newRow.N
    $ws.Range("O9").Value = $newRow.O
    $ws.Range("P9").Value = $newRow.P
    $ws.Range("Q9").Value = $newRow.Q
}
